$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 5 (pushes old rows 5-16 down to rows 10-21)
$ws.Rows.Item(5).Resize(5).Insert()

# Update existing rows 2-4 "Precio Ramo" (K column)
$ws.Range("K2").Value = 33000
$ws.Range("K3").Value = 29000
$ws.Range("K4").Value = 28000

# Fill new rows 5-9 with historical María González orders (PED016-PED020)
# Row 5: PED016
$ws.Range("A5").Value = 'PED016'
$ws.Range("B5").Value = '2025-07-10 11:00'
$ws.Range("C5").Value = '2025-07-11 15:00'
$ws.Range("D5").Value = 'JUEVES'
$ws.Range("E5").Value = 'WhatsApp'
$ws.Range("F5").Value = ''
$ws.Range("G5").Value = 'María González'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '+56912345678'
$ws.Range("H5").ClearFormats()
$ws.Range("I5").Value = 'Rosas Rojas Clásicas'
$ws.Range("J5").Value = 'Flores favoritas'
$ws.Range("K5").Value = 28000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 'Patricia González'
$ws.Range("N5").Value = 'Te extraño'
$ws.Range("O5").Value = 'María'
$ws.Range("P5").Value = 'Av. Apoquindo 1234, Las Condes'
$ws.Range("Q5").Value = 'Las Condes'
$ws.Range("R5").Value = 'Sin motivo'
$ws.Range("S5").Value = 'Archivado'
$ws.Range("T5").Value = 'Pagado'
$ws.Range("U5").Value = 'Normal'
$ws.Range("V5").Value = ''
$ws.Range("W5").Value = ''

# Row 6: PED017
$ws.Range("A6").Value = 'PED017'
$ws.Range("B6").Value = '2025-06-05 09:30'
$ws.Range("C6").Value = '2025-06-06 10:00'
$ws.Range("D6").Value = 'MIERCOLES'
$ws.Range("E6").Value = 'Shopify'
$ws.Range("F6").Value = '#SH1050'
$ws.Range("G6").Value = 'María González'
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '+56912345678'
$ws.Range("H6").ClearFormats()
$ws.Range("I6").Value = 'Jardín de Rosas'
$ws.Range("J6").Value = ''
$ws.Range("K6").Value = 28000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 'Elena González'
$ws.Range("N6").Value = 'Felicidades mamá'
$ws.Range("O6").Value = 'Tu hija'
$ws.Range("P6").Value = 'Av. Apoquindo 1234, Las Condes'
$ws.Range("Q6").Value = 'Las Condes'
$ws.Range("R6").Value = 'Día de la Madre'
$ws.Range("S6").Value = 'Despachados'
$ws.Range("T6").Value = 'Pagado'
$ws.Range("U6").Value = 'Normal'
$ws.Range("V6").Value = ''
$ws.Range("W6").Value = ''

# Row 7: PED018
$ws.Range("A7").Value = 'PED018'
$ws.Range("B7").Value = '2025-05-15 14:00'
$ws.Range("C7").Value = '2025-05-16 11:00'
$ws.Range("D7").Value = 'JUEVES'
$ws.Range("E7").Value = 'WhatsApp'
$ws.Range("F7").Value = ''
$ws.Range("G7").Value = 'María González'
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '+56912345678'
$ws.Range("H7").ClearFormats()
$ws.Range("I7").Value = 'Bouquet Romántico'
$ws.Range("J7").Value = 'Con tarjeta'
$ws.Range("K7").Value = 25000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 'Roberto González'
$ws.Range("N7").Value = 'Te amo'
$ws.Range("O7").Value = 'María'
$ws.Range("P7").Value = 'Av. Apoquindo 1234, Las Condes'
$ws.Range("Q7").Value = 'Las Condes'
$ws.Range("R7").Value = 'San Valentín'
$ws.Range("S7").Value = 'Archivado'
$ws.Range("T7").Value = 'Pagado'
$ws.Range("U7").Value = 'Normal'
$ws.Range("V7").Value = ''
$ws.Range("W7").Value = ''

# Row 8: PED019
$ws.Range("A8").Value = 'PED019'
$ws.Range("B8").Value = '2025-04-20 10:30'
$ws.Range("C8").Value = '2025-04-21 14:00'
$ws.Range("D8").Value = 'SABADO'
$ws.Range("E8").Value = 'Shopify'
$ws.Range("F8").Value = '#SH1020'
$ws.Range("G8").Value = 'María González'
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '+56912345678'
$ws.Range("H8").ClearFormats()
$ws.Range("I8").Value = 'Rosas Premium Mix'
$ws.Range("J8").Value = ''
$ws.Range("K8").Value = 28000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 'Sofía González'
$ws.Range("N8").Value = 'Feliz cumpleaños'
$ws.Range("O8").Value = 'María'
$ws.Range("P8").Value = 'Av. Apoquindo 1234, Las Condes'
$ws.Range("Q8").Value = 'Las Condes'
$ws.Range("R8").Value = 'Cumpleaños'
$ws.Range("S8").Value = 'Despachados'
$ws.Range("T8").Value = 'Pagado'
$ws.Range("U8").Value = 'Normal'
$ws.Range("V8").Value = ''
$ws.Range("W8").Value = ''

# Row 9: PED020
$ws.Range("A9").Value = 'PED020'
$ws.Range("B9").Value = '2025-03-12 16:00'
$ws.Range("C9").Value = '2025-03-13 09:00'
$ws.Range("D9").Value = 'MIERCOLES'
$ws.Range("E9").Value = 'WhatsApp'
$ws.Range("F9").Value = ''
$ws.Range("G9").Value = 'María González'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '+56912345678'
$ws.Range("H9").ClearFormats()
$ws.Range("I9").Value = 'Arreglo Elegante'
$ws.Range("J9").Value = 'Rosas rojas grandes'
$ws.Range("K9").Value = 25000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 'Andrés González'
$ws.Range("N9").Value = 'Gracias por todo'
$ws.Range("O9").Value = 'María'
$ws.Range("P9").Value = 'Av. Apoquindo 1234, Las Condes'
$ws.Range("Q9").Value = 'Las Condes'
$ws.Range("R9").Value = 'Agradecimiento'
$ws.Range("S9").Value = 'Archivado'
$ws.Range("T9").Value = 'Pagado'
$ws.Range("U9").Value = 'Normal'
$ws.Range("V9").Value = ''
$ws.Range("W9").Value = ''
